$wb = $excel.ActiveWorkbook

# --- Locate the existing "总计" (Total) sheet (3rd sheet, rId3 / sheetId 3) ---
$wsTotalOld = $wb.Worksheets.Item(3)

# --- Step 1: duplicate it to the end of the workbook BEFORE changing anything.
#     The duplicate inherits sheetPr / pageMargins / styles from the original
#     "总计" sheet, and will become the new "总计" sheet (with an extra
#     2022-Q1 summary row). The original sheet (3rd position) will be turned
#     into the new "2022-Q1" data sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsTotalOld.Copy($null, $lastSheet)
$wsTotalNew = $wb.Worksheets.Item($wb.Worksheets.Count)

# ===========================================================================
# Step 2: turn the original 3rd sheet into "2022-Q1" (fund holding detail).
#          Rename it *before* renaming the duplicate back to "总计" - both
#          sheets are briefly named "总计" after the Copy() above, so the
#          original must give up that name first to avoid a collision.
# ===========================================================================
$ws = $wsTotalOld
$ws.Name = "2022-Q1"
$wsTotalNew.Name = "总计"

# Copy the existing bold/bordered header style (from B1) onto the new header
# cells E1:H1 so they match the look of B1:D1.
$ws.Range("B1").Copy()
$ws.Range("E1:H1").PasteSpecial(-4122)

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Row 2: single fund-holding record. A2 already carries the index style (s=2)
# and value 0 from the original sheet, so it is left untouched.
$ws.Range("B2:G2").NumberFormat = "@"
$ws.Range("B2").Value = "002952"
$ws.Range("C2").Value = "建信多因子量化股票"
$ws.Range("D2").Value = "0.10"
$ws.Range("E2").Value = "91.47"
$ws.Range("F2").Value = "2.72"
$ws.Range("G2").Value = "0.0027"
$ws.Range("H2").Value = 10

# The old sheet had a 3rd row ("2021-Q2" totals) that no longer belongs here
# (new dimension is A1:H2) - remove it entirely.
$ws.Rows.Item(3).Delete()

# ===========================================================================
# Step 3: update the new "总计" sheet - insert a 2022-Q1 row at the top and
#          re-number the existing rows underneath it.
# ===========================================================================
$wt = $wsTotalNew
$wt.Rows.Item(2).Insert()

$wt.Range("A2").Value = 0
$wt.Range("B2").Value = "2022-Q1"
$wt.Range("C2").Value = 1
$wt.Range("D2").Value = 0

$wt.Range("A3").Value = 1
$wt.Range("B3").Value = "2021-Q3"
$wt.Range("C3").Value = 3
$wt.Range("D3").Value = 0.02

$wt.Range("A4").Value = 2
$wt.Range("B4").Value = "2021-Q2"
$wt.Range("C4").Value = 4
$wt.Range("D4").Value = 0.18
